$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiments")

# --- Row 36: new section header "NYU Depth v2 Labled Test Set (654 Entries)" ---
$ws.Range("A36").Value = "NYU Depth v2 Labled Test Set (654 Entries)"
$ws.Range("B36").Value = "delta1"
$ws.Range("C36").Value = "delta2"
$ws.Range("D36").Value = "delta3"
$ws.Range("E36").Value = "mse"
$ws.Range("F36").Value = "rmse"
$ws.Range("G36").Value = "rel_abs_dif"
$ws.Range("H36").Value = "rel_sqr_diff"
$ws.Range("I36").Value = "log10"
$ws.Range("J36").Value = "log_rmse"

# --- Row 38: DenseDepth_nohints results (filled in before row 37's note, matching authoring order) ---
$ws.Range("A38").Value = "DenseDepth_nohints"
$ws.Range("B38").Value = 0.85618639232578697
$ws.Range("C38").Value = 0.97836328187102095
$ws.Range("D38").Value = 0.99565340722963402
$ws.Range("E38").Value = 0.21353111323962001
$ws.Range("F38").Formula = "=SQRT(E38)"
$ws.Range("G38").Value = 0.11950074903631699
$ws.Range("H38").Value = 0.068724721778946202
$ws.Range("I38").Value = 0.051383407905070101
$ws.Range("J38").Value = 0.151374158361015
$ws.Range("L38").Value = "Uses rawdepth for evaluation (masks off invalid depth pixels), unlike what Wonka et.al. do in their paper."
$ws.Rows.Item(38).RowHeight = 51

# --- Row 37: DORN_nohints results run on the new labeled test set ---
$ws.Range("A37").Value = "DORN_nohints"
$ws.Range("B37").Value = 0.83929108669991404
$ws.Range("C37").Value = 0.95858464605035199
$ws.Range("D37").Value = 0.98564354107743601
$ws.Range("E37").Value = 0.26140680609123301
$ws.Range("F37").Formula = "=SQRT(E37)"
$ws.Range("G37").Value = 0.12950759483657001
$ws.Range("H37").Value = 0.087568223680133297
$ws.Range("I37").Value = 0.058503517225143298
$ws.Range("J37").Value = 0.17285743104927601
$ws.Range("L37").Value = "Pytorch version. Will run caffe version soon."
$ws.Rows.Item(37).RowHeight = 17

# --- View state: scrolled down a bit, selection on the newly added note cell ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("L38").Select()
